# Re-apply the "new" table style (PowerPoint's built-in "Medium Style 2 -
# Accent 1" GUID) to the three data tables that previously used the old
# default "Table_0" style GUID.  This mirrors the author re-styling those
# tables (slides 14, 15 and 16), each of which has the table as the first
# shape on the slide.

$p = $ppt.ActivePresentation

$oldStyleId = "{AE3AEF34-8078-402E-8BDC-628756838C81}"
$newStyleId = "{83FD9BE3-15D6-41B8-AC45-D24EA64D8B20}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $tbl = $shape.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}
